$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.528.45'
$ws.Range('E2').Value = '  +3.07%  '

$ws.Range('D3').Value = '1.818.72'
$ws.Range('E3').Value = '  +4.24%  '

$ws.Range('D4').Value = '1.002'

$ws.Range('D5').Value = '343.76'
$ws.Range('E5').Value = '  +2.80%  '

$ws.Range('E6').Value = '  +0.49%  '

$ws.Range('D7').Value = '0.3833'
$ws.Range('E7').Value = '  +2.83%  '

$ws.Range('D8').Value = '0.3536'
$ws.Range('E8').Value = '  +3.79%  '

$ws.Range('D9').Value = '48.99'
$ws.Range('E9').Value = '  -1.62%  '

$ws.Range('D10').Value = '1.237'
$ws.Range('E10').Value = '  +2.25%  '

$ws.Range('D11').Value = '0.07772'
$ws.Range('E11').Value = '  +2.97%  '

$ws.Range('E12').Value = '  +0.37%  '

$ws.Range('D13').Value = '22.34'
$ws.Range('E13').Value = '  +8.38%  '

$ws.Range('D14').Value = '6.598'
$ws.Range('E14').Value = '  +2.12%  '

$ws.Range('D15').Value = '1.816.84'
$ws.Range('E15').Value = '  +4.59%  '

$ws.Range('D16').Value = '7.209'
$ws.Range('E16').Value = '  +2.11%  '

$ws.Range('D17').Value = '0.00001123'

$ws.Range('D18').Value = '0.06728'
$ws.Range('E18').Value = '  +0.34%  '

$ws.Range('D19').Value = '86.57'
$ws.Range('E19').Value = '  +2.79%  '

$ws.Range('D20').Value = '0.9996'
$ws.Range('E20').Value = '  +0.41%  '

$ws.Range('E21').Value = '  +5.16%  '

$ws.Range('D22').Value = '6.553'
$ws.Range('E22').Value = '  +5.47%  '

$ws.Range('D23').Value = '13.19'
$ws.Range('E23').Value = '  +0.38%  '

$ws.Range('D24').Value = '27.509.82'
$ws.Range('E24').Value = '  +3.31%  '

$ws.Range('D25').Value = '2.463'
$ws.Range('E25').Value = '  -0.47%  '

$ws.Range('D26').Value = '2.690'
$ws.Range('E26').Value = '  +6.24%  '

$ws.Range('D27').Value = '22.21'
$ws.Range('E27').Value = '  +12.33%  '

$ws.Range('E28').Value = '  +4.02%  '

$ws.Range('D29').Value = '154.14'
$ws.Range('E29').Value = '  +1.37%  '

$ws.Range('D30').Value = '2.023.31'
$ws.Range('E30').Value = '  +4.81%  '

$ws.Range('D31').Value = '136.11'
$ws.Range('E31').Value = '  +2.76%  '

$ws.Range('D32').Value = '6.382'
$ws.Range('E32').Value = '  +1.63%  '

$ws.Range('D33').Value = '4.077'
$ws.Range('E33').Value = '  -1.19%  '

$ws.Range('D34').Value = '13.92'
$ws.Range('E34').Value = '  +5.18%  '

$ws.Range('D35').Value = '0.08812'
$ws.Range('E35').Value = '  +2.73%  '

$ws.Range('D36').Value = '1.689'
$ws.Range('E36').Value = '  -1.17%  '

$ws.Range('D37').Value = '5.637'
$ws.Range('E37').Value = '  +2.40%  '

$ws.Range('D38').Value = '0.7037'
$ws.Range('E38').Value = '  +11.93%  '

$ws.Range('D39').Value = '0.2264'
$ws.Range('E39').Value = '  +3.67%  '

$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.06514'
$ws.Range('E40').Value = '  +2.11%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.02401'
$ws.Range('E41').Value = '  +1.72%  '

$ws.Range('D42').Value = '8.975'
$ws.Range('E42').Value = '  +3.10%  '

$ws.Range('D43').Value = '1.299'
$ws.Range('E43').Value = '  +4.35%  '

$ws.Range('D44').Value = '14.77'
$ws.Range('E44').Value = '  +1.54%  '

$ws.Range('D45').Value = '0.6618'
$ws.Range('E45').Value = '  +7.99%  '

$ws.Range('D46').Value = '0.9995'
$ws.Range('E46').Value = '  +0.37%  '

$ws.Range('D47').Value = '3.961'
$ws.Range('E47').Value = '  +1.49%  '

$ws.Range('D48').Value = '2.192'
$ws.Range('E48').Value = '  +5.68%  '

$ws.Range('D49').Value = '132.61'
$ws.Range('E49').Value = '  +2.42%  '

$ws.Range('D50').Value = '0.07339'
$ws.Range('E50').Value = '  -0.20%  '

$ws.Range('D51').Value = '80.77'
$ws.Range('E51').Value = '  +3.35%  '
